$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 18 de Agosto de 2020 a las 07:56"

# Re-rank Butan / Brunei (rows 190-191) - Butan overtook Brunei with updated case counts
$ws.Range("A190").Value = "Butan"
$ws.Range("B190").Value = 146
$ws.Range("C190").Value = 5
$ws.Range("D190").Value = 103
$ws.Range("E190").Value = 43
$ws.Range("H190").Value = 0

$ws.Range("A191").Value = "Brunei"
$ws.Range("B191").Value = 142
$ws.Range("D191").Value = 139
$ws.Range("E191").Value = 0
$ws.Range("H191").Value = 3

# Re-rank Islas Malvinas / Montserrat (rows 213-214)
$ws.Range("A213").Value = "Islas Malvinas"
$ws.Range("D213").Value = 13
$ws.Range("H213").Value = 0

$ws.Range("A214").Value = "Montserrat"
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1

# Updated case counts for several countries (rows 6, 18, 33, 61, 73)
$ws.Range("B6").Value = 2703517
$ws.Range("C6").Value = 1913
$ws.Range("D6").Value = 1977779
$ws.Range("E6").Value = 673802
$ws.Range("G6").Value = 11
$ws.Range("H6").Value = 51936

$ws.Range("B18").Value = 289832
$ws.Range("C18").Value = 617
$ws.Range("D18").Value = 270009
$ws.Range("E18").Value = 13633
$ws.Range("G18").Value = 15
$ws.Range("H18").Value = 6190

$ws.Range("B33").Value = 95129
$ws.Range("C33").Value = 378
$ws.Range("D33").Value = 71165
$ws.Range("E33").Value = 23272

$ws.Range("B61").Value = 36100
$ws.Range("C61").Value = 398
$ws.Range("D61").Value = 31580
$ws.Range("E61").Value = 4280
$ws.Range("G61").Value = 4
$ws.Range("H61").Value = 240

$ws.Range("D73").Value = 11017
$ws.Range("E73").Value = 11551
$ws.Range("G73").Value = 7
$ws.Range("H73").Value = 625
